# TC28_Verify_PlaceOrder_SelfService_SingleUser.xlsx
# "Changes done for remaining KIT testcases"
#
# The test-step renamed the web-element/object used to verify a successful
# order placement (PlaceOrderMsg/PlaceOrderSuccessNew -> PlaceOrderSuccess,
# and its message key PlaceOrderMsg -> PlaceOrderSuccessMSG / text
# "Your Order has been Placed!" -> "Thank You"), merging what used to be
# three verification rows (53-55) into two (53-54) on the first sheet, and
# updating the corresponding lookup row on the Testdata sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC28_Verify_PlaceOrder_SS_SU")
$ws2 = $wb.Worksheets.Item("Testdata")

# --- Sheet 1: TC28_Verify_PlaceOrder_SS_SU -------------------------------
# Row 53 used to verify presence of the "PlaceOrderMsg" element; now it
# verifies the "PlaceOrderSuccess" element and reports the
# "PlaceOrderSuccessMSG" data descriptor.
$ws1.Range("C53").Value = "PlaceOrderSuccess"
$ws1.Range("E53").Value = "PlaceOrderSuccessMSG"

# Row 54 used to be a second VERIFY_WEBELEMENT_PRESENT check on
# "PlaceOrderSuccessNew"; it is now a VERIFY_TEXT_PRESENT check on
# "PlaceOrderSuccess" with the "PlaceOrderSuccessMSG" descriptor.
$ws1.Range("B54").Value = "VERIFY_TEXT_PRESENT"
$ws1.Range("C54").Value = "PlaceOrderSuccess"
$ws1.Range("E54").Value = "PlaceOrderSuccessMSG"

# Old row 55 (VERIFY_TEXT_PRESENT / PlaceOrderMsg / CSS / PlaceOrderMsg) is
# now redundant with the merged row 54 above, so it is removed and the
# rows below it (none, in this sheet) shift up.
$ws1.Rows.Item(55).Delete()

$ws1.Range("E54").Select()

# --- Sheet 2: Testdata ----------------------------------------------------
# The lookup row for the success message is renamed/re-worded to match.
$ws2.Range("A57").Value = "PlaceOrderSuccessMSG"
$ws2.Range("B57").Value = "Thank You"

$ws2.Range("A57").Select()
$ws1.Activate()
